$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data: date 2012-07-31 and activity description
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A7").Value = 41121
$ws.Range("B7").Value = "Implemented Bealto ParallelSelectionSortLocal and ParallelSelectionSortBlocks"

# Update the active selection on the sheet
$ws.Range("B18").Select()
